$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "299.84") need the
# column pre-set to Text format, otherwise Excel auto-converts the literal
# into a numeric value instead of keeping it as the original text content.
$textCells = @("D5","D6","D7","D9","D10","D11","D12","D14","D16","D18","D21","D22","D23","D24","D26","D28","D30","D31","D32","D33","D34","D36","D37","D38","D39","D40","D41","D45","D46","D47","D48","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '42.390.02'
$ws.Range('E2').Value = '  -2.69%  '
$ws.Range('D3').Value = '2.274.84'
$ws.Range('E3').Value = '  -4.49%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '299.84'
$ws.Range('E5').Value = '  -3.26%  '
$ws.Range('D6').Value = '96.75'
$ws.Range('E6').Value = '  -7.50%  '
$ws.Range('D7').Value = '0.503'
$ws.Range('E7').Value = '  -1.60%  '
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').Value = '  -5.10%  '
$ws.Range('D10').Value = '33.60'
$ws.Range('E10').Value = '  -6.67%  '
$ws.Range('D11').Value = '50.57'
$ws.Range('E11').Value = '  -5.34%  '
$ws.Range('D12').Value = '0.0786'
$ws.Range('E12').Value = '  -3.34%  '
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').Value = '6.63'
$ws.Range('E14').Value = '  -4.99%  '
$ws.Range('D15').Value = '2.627.79'
$ws.Range('E15').Value = '  -4.52%  '
$ws.Range('D16').Value = '15.16'
$ws.Range('E16').Value = '  -2.85%  '
$ws.Range('D17').Value = '2.272.07'
$ws.Range('E17').Value = '  -4.83%  '
$ws.Range('D18').Value = '0.784'
$ws.Range('E18').Value = '  -3.41%  '
$ws.Range('D19').Value = '42.322.31'
$ws.Range('E19').Value = '  -2.83%  '
$ws.Range('D20').Value = '0.0₃0890'
$ws.Range('E20').Value = '  -2.91%  '
$ws.Range('D21').Value = '11.38'
$ws.Range('E21').Value = '  -4.42%  '
$ws.Range('D22').Value = '5.98'
$ws.Range('E22').Value = '  -5.31%  '
$ws.Range('D23').Value = '66.42'
$ws.Range('E23').Value = '  -2.90%  '
$ws.Range('D24').Value = '234.49'
$ws.Range('E24').Value = '  -2.56%  '
$ws.Range('E25').Value = '  -6.03%  '
$ws.Range('D26').Value = '2.48'
$ws.Range('E26').Value = '  -5.16%  '
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').Value = '24.29'
$ws.Range('E28').Value = '  -6.06%  '
$ws.Range('D30').Value = '33.65'
$ws.Range('E30').Value = '  -7.98%  '
$ws.Range('D31').Value = '163.40'
$ws.Range('E31').Value = '  +1.62%  '
$ws.Range('D32').Value = '9.06'
$ws.Range('E32').Value = '  -4.85%  '
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').Value = '4.93'
$ws.Range('E34').Value = '  -6.18%  '
$ws.Range('E35').Value = '  -4.82%  '
$ws.Range('D36').Value = '0.0694'
$ws.Range('E36').Value = '  -5.96%  '
$ws.Range('D37').Value = '4.34'
$ws.Range('E37').Value = '  -7.67%  '
$ws.Range('D38').Value = '2.82'
$ws.Range('E38').Value = '  -9.61%  '
$ws.Range('D39').Value = '16.01'
$ws.Range('E39').Value = '  -12.74%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '0.0997'
$ws.Range('E40').Value = '  -5.92%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '1.76'
$ws.Range('E41').Value = '  -9.31%  '
$ws.Range('E42').Value = '  -3.68%  '
$ws.Range('E43').Value = '  -7.65%  '
$ws.Range('D44').Value = '1.958.28'
$ws.Range('E44').Value = '  -3.82%  '
$ws.Range('D45').Value = '0.0280'
$ws.Range('E45').Value = '  -3.44%  '
$ws.Range('D46').Value = '17.71'
$ws.Range('E46').Value = '  -10.37%  '
$ws.Range('D47').Value = '9.63'
$ws.Range('E47').Value = '  -8.92%  '
$ws.Range('D48').Value = '2.80'
$ws.Range('E48').Value = '  -10.32%  '
$ws.Range('E49').Value = '  -4.71%  '
$ws.Range('D50').Value = '4.67'
$ws.Range('E50').Value = '  -1.53%  '
$ws.Range('D51').Value = '2.500.57'
